# Update scores on the "ScoreF" worksheet and adjust the active tab/selection
# so they match the committed workbook state.

$wb = $excel.ActiveWorkbook

$wsF = $wb.Worksheets.Item("ScoreF")

# N/O column corrections and new P (total) column values for ScoreF.
$wsF.Range("N3").Value = 5
$wsF.Range("O3").Value = 50

$wsF.Range("N9").Value = 6

$wsF.Range("N17").Value = 7
$wsF.Range("O17").Value = 54

$wsF.Range("N19").Value = 5
$wsF.Range("O19").Value = 35

$wsF.Range("N21").Value = 7
$wsF.Range("O21").Value = 4

$wsF.Range("N23").Value = 7
$wsF.Range("O23").Value = 34

$wsF.Range("P2").Value = 120
$wsF.Range("P3").Value = 156
$wsF.Range("P4").Value = 136
$wsF.Range("P5").Value = 149
$wsF.Range("P6").Value = 130
$wsF.Range("P7").Value = 138
$wsF.Range("P8").Value = 144
$wsF.Range("P9").Value = 156
$wsF.Range("P10").Value = 136
$wsF.Range("P11").Value = 143
$wsF.Range("P12").Value = 136
$wsF.Range("P13").Value = 152
$wsF.Range("P14").Value = 152
$wsF.Range("P15").Value = 151
$wsF.Range("P16").Value = 136
$wsF.Range("P17").Value = 156
$wsF.Range("P18").Value = 130
$wsF.Range("P19").Value = 156
$wsF.Range("P20").Value = 151
$wsF.Range("P21").Value = 156
$wsF.Range("P22").Value = 136
$wsF.Range("P23").Value = 156
$wsF.Range("P24").Value = 144
$wsF.Range("P25").Value = 140

# Move the active tab/selection from "ScoreM" to "ScoreF" (the last-saved view).
$wsF.Range("P20").Select()
$wsF.Activate()
